$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove sheet protection (password) so the new column/content can be added;
# the protection flags themselves (formatCells/formatColumns/etc.) stay "unprotected".
$ws.Unprotect()

# Add the new "Supervisor" header in column E (new shared string value)
$ws.Range("E1").Value = "Supervisor"

# Apply the same header style (fill/border/bold font) used by the other header cells
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the width for the new column E (target stored width ~85.7109375 chars)
$ws.Range("E1").EntireColumn.ColumnWidth = 84.8776

# Add a comment to E1 describing how to fill in supervisors
$comment = $ws.Range("E1").AddComment("Thomas Ingeman-Nielsen:" + [char]10 + [char]10 + "List of supervisors separated by commas:" + [char]10 + "First Middle Last, First Middle Last, ...")

# Update selection to match the target state
$ws.Range("C3").Select()

# Re-apply sheet protection without a password (matches diff: protection flags kept, password removed)
$ws.Protect($null, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
